# Tugas UTS-Agile.docx edit
#
# 1. Fix the name at the end of the document: "Alwi nt" -> "Alwi ngga nt"
# 2. Make the page orientation explicit (portrait) on the section's page setup.

$d = $word.ActiveDocument

# --- 1. Correct the trailing name run ("Alwi nt" -> "Alwi ngga nt") ---
$d.Content.Find.ClearFormatting()
$found = $d.Content.Find.Execute(
    "Alwi nt",   # FindText
    $false,      # MatchCase
    $true,       # MatchWholeWord
    $false,      # MatchWildcards
    $false,      # MatchSoundsLike
    $false,      # MatchAllWordForms
    $true,       # Forward
    1,           # Wrap (wdFindContinue)
    $false,      # Format
    "Alwi ngga nt", # ReplaceWith
    2            # Replace (wdReplaceAll)
)

# --- 2. Make the (already-portrait) page orientation explicit in the XML ---
$d.PageSetup.Orientation = 0   # wdOrientPortrait
